$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "id_produto" column (column C). This shifts the remaining
# columns one place to the left: Categoria D->C, População E->D,
# unidade F->E, ano G->F.
$ws.Range("C1:C28").Delete()

# After the shift, the "População" numbers (now in column D, rows 2-28)
# must be stored as text rather than numeric values.
$rng = $ws.Range("D2:D28")
$rng.NumberFormat = "@"
foreach ($cell in $rng) {
    $cell.Value = $cell.Text
}
$rng.ClearFormats()
